$wb = $excel.ActiveWorkbook

# 1. Rename the "APPL" sheet to "AAPL" (typo fix)
$wsAAPL = $wb.Worksheets.Item("APPL")
$wsAAPL.Name = "AAPL"

# 2. Fix stray trailing-space duplicate of the "Rating" header on the PG sheet
#    (was "Rating ", now matches the "Rating" text used on every other sheet).
$wsPG = $wb.Worksheets.Item("PG")
$wsPG.Range("D1").Value = "Rating"

# 3. Update the remembered cell selection on a handful of sheets to D1
[void]$wsAAPL.Range("D1").Select()

$wsMSFT = $wb.Worksheets.Item("MSFT")
[void]$wsMSFT.Range("D1").Select()

$wsIBM = $wb.Worksheets.Item("IBM")
[void]$wsIBM.Range("D1").Select()

[void]$wsPG.Range("D1").Select()

# Restore the originally active sheet/tab
$wsAll = $wb.Worksheets.Item("All")
[void]$wsAll.Activate()
